# Auto-generated edit script applying the Anima_Profits.xlsx market-price update diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 400
$ws.Range("I18").Value = 300
$ws.Range("K18").Value = 300
$ws.Range("M18").Value = -16
$ws.Range("H33").Value = 1351.5625
$ws.Range("I33").Value = 871
$ws.Range("J33").Value = 3434
$ws.Range("K33").Value = 871
$ws.Range("L33").Value = 3434
$ws.Range("M33").Value = -642
$ws.Range("N33").Value = -3892
$ws.Range("H109").Value = 20000
$ws.Range("J109").Value = 20000
$ws.Range("L109").Value = 20000
$ws.Range("N109").Value = -22774
$ws.Range("H116").Value = 2670
$ws.Range("I116").Value = 2342.5
$ws.Range("J116").Value = 3325
$ws.Range("K116").Value = 2342.5
$ws.Range("L116").Value = 3325
$ws.Range("M116").Value = 1099.5
$ws.Range("N116").Value = -10209
$ws.Range("H127").Value = 1229.8889
$ws.Range("I127").Value = 733.1667
$ws.Range("J127").Value = 1478.25
$ws.Range("K127").Value = 2199.5001
$ws.Range("L127").Value = 4434.75
$ws.Range("M127").Value = 2760.4999
$ws.Range("N127").Value = -14354.75
$ws.Range("H137").Value = 3027.4473
$ws.Range("I137").Value = 2260.7666
$ws.Range("J137").Value = 5902.5
$ws.Range("K137").Value = 6782.2998
$ws.Range("L137").Value = 17707.5
$ws.Range("M137").Value = -4232.2998
$ws.Range("N137").Value = -22807.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1171.0322
$ws.Range("I74").Value = 808.087
$ws.Range("J74").Value = 2214.5
$ws.Range("K74").Value = 808.087
$ws.Range("L74").Value = 2214.5
$ws.Range("M74").Value = 65.91300000000001
$ws.Range("N74").Value = -3962.5
$ws.Range("H77").Value = 1171.0322
$ws.Range("I77").Value = 808.087
$ws.Range("J77").Value = 2214.5
$ws.Range("K77").Value = 4040.435
$ws.Range("L77").Value = 11072.5
$ws.Range("M77").Value = 327.5650000000001
$ws.Range("N77").Value = -19808.5
$ws.Range("H132").Value = 3434.9138
$ws.Range("I132").Value = 2470.0698
$ws.Range("K132").Value = 7410.209400000001
$ws.Range("M132").Value = -4880.209400000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H93").Value = 64925
$ws.Range("J93").Value = 64925
$ws.Range("L93").Value = 64925
$ws.Range("N93").Value = -68669
$ws.Range("H97").Value = 20000
$ws.Range("I97").Value = 20000
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 20000
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -19009
$ws.Range("N97").ClearContents()
$ws.Range("H99").Value = 1545.091
$ws.Range("I99").Value = 1549.6
$ws.Range("K99").Value = 1549.6
$ws.Range("M99").Value = -51.59999999999991
$ws.Range("H109").Value = 40349.92
$ws.Range("J109").Value = 40349.92
$ws.Range("L109").Value = 40349.92
$ws.Range("N109").Value = -43123.92
$ws.Range("H123").Value = 26490
$ws.Range("J123").Value = 26490
$ws.Range("L123").Value = 26490
$ws.Range("N123").Value = -36290

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7024.8296
$ws.Range("I31").Value = 1309.0869
$ws.Range("K31").Value = 1309.0869
$ws.Range("M31").Value = -1014.0869
$ws.Range("H34").Value = 7024.8296
$ws.Range("I34").Value = 1309.0869
$ws.Range("K34").Value = 1309.0869
$ws.Range("M34").Value = -1107.0869

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H39").Value = 2036.1818
$ws.Range("J39").Value = 2036.1818
$ws.Range("L39").Value = 6108.5454
$ws.Range("N39").Value = -6696.5454
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H92").Value = 585.5
$ws.Range("J92").Value = 575
$ws.Range("L92").Value = 1725
$ws.Range("N92").Value = -4221
$ws.Range("H110").Value = 13296.52
$ws.Range("J110").Value = 13358.917
$ws.Range("L110").Value = 40076.751
$ws.Range("N110").Value = -48256.751
$ws.Range("H131").Value = 1050.7778
$ws.Range("J131").Value = 1062.3922
$ws.Range("L131").Value = 3187.1766
$ws.Range("N131").Value = -13267.1766
$ws.Range("H139").Value = 3362.4524
$ws.Range("I139").Value = 1430
$ws.Range("J139").Value = 4228.724
$ws.Range("K139").Value = 4290
$ws.Range("L139").Value = 12686.172
$ws.Range("M139").Value = 850
$ws.Range("N139").Value = -22966.172
$ws.Range("H141").Value = 5127.647
$ws.Range("I141").Value = 2117
$ws.Range("K141").Value = 6351
$ws.Range("M141").Value = -1171

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 48778.285
$ws.Range("J94").Value = 48778.285
$ws.Range("L94").Value = 48778.285
$ws.Range("N94").Value = -50130.285
$ws.Range("H112").Value = 43632.668
$ws.Range("J112").Value = 43632.668
$ws.Range("L112").Value = 43632.668
$ws.Range("N112").Value = -45848.668
$ws.Range("H118").Value = 74810
$ws.Range("J118").Value = 74810
$ws.Range("L118").Value = 74810
$ws.Range("N118").Value = -78124
$ws.Range("H124").Value = 47896
$ws.Range("J124").Value = 47896
$ws.Range("L124").Value = 47896
$ws.Range("N124").Value = -57716
$ws.Range("H132").Value = 2069.524
$ws.Range("I132").Value = 1716.9333
$ws.Range("J132").Value = 2951
$ws.Range("K132").Value = 5150.7999
$ws.Range("L132").Value = 8853
$ws.Range("M132").Value = -2620.7999
$ws.Range("N132").Value = -13913
$ws.Range("H133").Value = 50780
$ws.Range("J133").Value = 50780
$ws.Range("L133").Value = 50780
$ws.Range("N133").Value = -60900

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 960
$ws.Range("I46").Value = 700
$ws.Range("J46").Value = 1025
$ws.Range("K46").Value = 700
$ws.Range("L46").Value = 1025
$ws.Range("M46").Value = -512
$ws.Range("N46").Value = -1401
$ws.Range("H122").Value = 3451.75
$ws.Range("I122").Value = 3213.7144
$ws.Range("J122").Value = 3785
$ws.Range("K122").Value = 9641.143199999999
$ws.Range("L122").Value = 11355
$ws.Range("M122").Value = -7191.143199999999
$ws.Range("N122").Value = -16255

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 861.28
$ws.Range("I113").Value = 1048.5555
$ws.Range("K113").Value = 3145.6665
$ws.Range("M113").Value = -975.6664999999998
